$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 477; this shifts the existing rows 477-545
# down to 478-546 (Excel copies formatting from the row above by default).
$ws.Rows.Item(477).Insert()

# Populate the newly inserted row 477 with the new data record.
$ws.Range("A477").Value = 10
$ws.Range("B477").Value = "Vega Modelo de Temuco"
$ws.Range("C477").Value = "La Araucanía"
$ws.Range("D477").Value = 45034
$ws.Range("E477").Value = 9
$ws.Range("F477").Value = 100112040
$ws.Range("G477").Value = "Cilantro"
$ws.Range("H477").Value = "Sin especificar"
$ws.Range("I477").Value = "Primera"
$ws.Range("J477").Value = 85
$ws.Range("K477").Value = 5000
$ws.Range("L477").Value = 5000
$ws.Range("M477").Value = 5000
$ws.Range("N477").Value = "$/docena de atados (2 kilos)"
$ws.Range("O477").Value = "Provincia de Cautín"
$ws.Range("P477").Value = 2500
$ws.Range("Q477").Value = 2
$ws.Range("R477").Value = "Hortaliza"
